# Update MapItem B column values (Id) to the new numbering scheme
# and restore the active cell selection to B21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5; Value = 100101 },
    @{ Row = 6; Value = 100102 },
    @{ Row = 7; Value = 100103 },
    @{ Row = 8; Value = 100104 },
    @{ Row = 9; Value = 100105 },
    @{ Row = 10; Value = 100201 },
    @{ Row = 11; Value = 100202 },
    @{ Row = 12; Value = 100203 },
    @{ Row = 13; Value = 100301 },
    @{ Row = 14; Value = 100302 },
    @{ Row = 15; Value = 100303 },
    @{ Row = 16; Value = 100304 },
    @{ Row = 17; Value = 100305 },
    @{ Row = 18; Value = 100306 },
    @{ Row = 19; Value = 100307 },
    @{ Row = 20; Value = 100308 },
    @{ Row = 21; Value = 100309 },
    @{ Row = 22; Value = 100401 },
    @{ Row = 23; Value = 100402 },
    @{ Row = 24; Value = 100403 },
    @{ Row = 25; Value = 100404 },
    @{ Row = 26; Value = 100405 },
    @{ Row = 27; Value = 100406 },
    @{ Row = 28; Value = 100407 },
    @{ Row = 29; Value = 100408 },
    @{ Row = 30; Value = 100409 },
    @{ Row = 31; Value = 100410 },
    @{ Row = 32; Value = 100501 },
    @{ Row = 33; Value = 100502 },
    @{ Row = 34; Value = 100503 },
    @{ Row = 35; Value = 100504 },
    @{ Row = 36; Value = 100505 },
    @{ Row = 37; Value = 100601 },
    @{ Row = 38; Value = 100602 },
    @{ Row = 39; Value = 100601 },
    @{ Row = 40; Value = 100602 },
    @{ Row = 41; Value = 100603 },
    @{ Row = 42; Value = 100604 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.Value
}

$ws.Range("B21").Select()
